# Aggiornamento fino a 02/05
# Append 6 new daily rows (239-244) to the data sheet, extending the
# used range from A1:D238 to A1:D244, mirroring the style of the
# preceding rows (date column uses style index "2" already applied
# via the existing formatted cells above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row from the existing data (row 238) and copy its
# date-cell formatting down for the new rows so new cells in column A
# keep the same style as the rows above them.
$lastRow = 238

$newData = @(
    @{ Row = 239; A = 44313; B = 0; C = 2; D = 57.75339301183945 },
    @{ Row = 240; A = 44314; B = 0; C = 1; D = 28.87669650591972 },
    @{ Row = 241; A = 44315; B = 0; C = 1; D = 28.87669650591972 },
    @{ Row = 242; A = 44316; B = 0; C = 0; D = 0 },
    @{ Row = 243; A = 44317; B = 1; C = 1; D = 28.87669650591972 },
    @{ Row = 244; A = 44318; B = 0; C = 1; D = 28.87669650591972 }
)

foreach ($entry in $newData) {
    $r = $entry.Row

    # Copy the formatting of the previous row's cells down to the new
    # row (this preserves the bold/centered/bordered date style used
    # by column A, and the plain style of columns B-D).
    $ws.Range("A$lastRow`:D$lastRow").Copy() | Out-Null
    $ws.Range("A$r`:D$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D

    $lastRow = $r
}

$excel.CutCopyMode = $false
